$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("J3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

$ws.Range("D4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("J4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

$ws.Range("D5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("J5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

$ws.Range("D6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("J6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

$ws.Range("D8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("J8").Value = "15,23 TL - 30,47 TL - 304,72 TL"

$ws.Range("D9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("J9").Value = "15,23 TL - 30,47 TL - 304,72 TL"

$ws.Range("D10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("J10").Value = "15,23 TL - 30,47 TL - 304,72 TL"

$ws.Range("D11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("J11").Value = "3,05 TL - 6,1 TL - 76,18 TL"

$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"

$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("J13").Value = "Hesaba: Asgari 1 TL | Azami 995,5 TL"

$ws.Range("D14").Value = "2.300 TL - 9.500 TL"
$ws.Range("J14").Value = "1.188 TL - 593 TL"
